$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price column cells that are being updated,
# so numeric-looking values (e.g. "1.003") are not auto-converted to numbers.
$priceCells = @("D2", "D3", "D4", "D5", "D7", "D8", "D9", "D10", "D11", "D14", "D15", "D16", "D17", "D18", "D19", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D44", "D45", "D46", "D47", "D48", "D49", "D51")
foreach ($cell in $priceCells) {
    $ws.Range($cell).NumberFormat = "@"
}

# Apply the updated cell values from the source diff.
$ws.Range('D2').Value = '27.534.33'
$ws.Range('D3').Value = '1.751.50'
$ws.Range('E3').Value = '  -3.49%  '
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').Value = '324.23'
$ws.Range('E5').Value = '  -0.50%  '
$ws.Range('E6').Value = '  +0.14%  '
$ws.Range('D7').Value = '0.4473'
$ws.Range('E7').Value = '  +2.12%  '
$ws.Range('D8').Value = '0.3614'
$ws.Range('E8').Value = '  -1.54%  '
$ws.Range('D9').Value = '0.07493'
$ws.Range('E9').Value = '  -2.34%  '
$ws.Range('D10').Value = '42.13'
$ws.Range('E10').Value = '  -5.89%  '
$ws.Range('D11').Value = '1.103'
$ws.Range('E11').Value = '  -3.18%  '
$ws.Range('E12').Value = '  +0.08%  '
$ws.Range('E13').Value = '  -5.92%  '
$ws.Range('D14').Value = '6.050'
$ws.Range('E14').Value = '  -4.08%  '
$ws.Range('D15').Value = '7.185'
$ws.Range('E15').Value = '  -4.03%  '
$ws.Range('D16').Value = '1.750.27'
$ws.Range('E16').Value = '  -3.84%  '
$ws.Range('D17').Value = '92.73'
$ws.Range('E17').Value = '  -2.43%  '
$ws.Range('D18').Value = '0.00001063'
$ws.Range('E18').Value = '  -1.50%  '
$ws.Range('D19').Value = '0.06414'
$ws.Range('E19').Value = '  -0.84%  '
$ws.Range('E20').Value = '  +0.11%  '
$ws.Range('D21').Value = '17.04'
$ws.Range('E21').Value = '  -1.89%  '
$ws.Range('D22').Value = '5.855'
$ws.Range('E22').Value = '  -6.19%  '
$ws.Range('D23').Value = '27.590.41'
$ws.Range('E23').Value = '  -2.48%  '
$ws.Range('D24').Value = '11.23'
$ws.Range('E24').Value = '  -2.84%  '
$ws.Range('D25').Value = '2.103'
$ws.Range('E25').Value = '  -1.02%  '
$ws.Range('D26').Value = '162.13'
$ws.Range('E26').Value = '  +0.68%  '
$ws.Range('D27').Value = '20.44'
$ws.Range('E27').Value = '  -1.31%  '
$ws.Range('D28').Value = '1.951.61'
$ws.Range('E28').Value = '  -3.73%  '
$ws.Range('D29').Value = '2.128'
$ws.Range('E29').Value = '  -6.45%  '
$ws.Range('D30').Value = '125.21'
$ws.Range('E30').Value = '  -3.26%  '
$ws.Range('D31').Value = '1.081'
$ws.Range('E31').Value = '  -10.85%  '
$ws.Range('D32').Value = '0.09032'
$ws.Range('E32').Value = '  -1.18%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').Value = '5.554'
$ws.Range('E33').Value = '  -7.33%  '
$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D34').Value = '3.635'
$ws.Range('E34').Value = '  +2.72%  '
$ws.Range('D35').Value = '12.01'
$ws.Range('E35').Value = '  -7.76%  '
$ws.Range('D36').Value = '0.02313'
$ws.Range('E36').Value = '  -2.08%  '
$ws.Range('D37').Value = '0.2094'
$ws.Range('E37').Value = '  -3.59%  '
$ws.Range('D38').Value = '0.6380'
$ws.Range('E38').Value = '  -3.35%  '
$ws.Range('D39').Value = '0.05979'
$ws.Range('E39').Value = '  -3.68%  '
$ws.Range('B40').Value = 'InternetComputer(DFINITY)'
$ws.Range('C40').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D40').Value = '4.961'
$ws.Range('E40').Value = '  -5.36%  '
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').Value = '1.208'
$ws.Range('E41').Value = '  +0.74%  '
$ws.Range('D42').Value = '1.000'
$ws.Range('E42').Value = '  +0.11%  '
$ws.Range('E43').Value = '  -3.02%  '
$ws.Range('D44').Value = '7.795'
$ws.Range('E44').Value = '  -3.27%  '
$ws.Range('D45').Value = '13.19'
$ws.Range('E45').Value = '  -4.82%  '
$ws.Range('B46').Value = 'PancakeSwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D46').Value = '3.713'
$ws.Range('E46').Value = '  -0.60%  '
$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').Value = '0.5876'
$ws.Range('E47').Value = '  -3.67%  '
$ws.Range('D48').Value = '1.956'
$ws.Range('E48').Value = '  -3.19%  '
$ws.Range('D49').Value = '121.37'
$ws.Range('E49').Value = '  -3.14%  '
$ws.Range('E50').Value = '  +0.11%  '
$ws.Range('D51').Value = '0.06872'
$ws.Range('E51').Value = '  -1.73%  '
